$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1314.5883
$ws.Range("I53").Value = 915.61536
$ws.Range("J53").Value = 2611.25
$ws.Range("K53").Value = 915.61536
$ws.Range("L53").Value = 2611.25
$ws.Range("M53").Value = -278.61536
$ws.Range("N53").Value = -3885.25

$ws.Range("H62").Value = 2416
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2499.2
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2499.2
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3747.2

$ws.Range("H65").Value = 2416
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2499.2
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 12496
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -18736

$ws.Range("H111").Value = 3883.8
$ws.Range("I111").Value = 5273
$ws.Range("J111").Value = 1800
$ws.Range("K111").Value = 15819
$ws.Range("L111").Value = 5400
$ws.Range("M111").Value = -12752
$ws.Range("N111").Value = -11534

$ws.Range("H132").Value = 7411077.5
$ws.Range("I132").Value = 7938602
$ws.Range("K132").Value = 23815806
$ws.Range("M132").Value = -23813276

$ws.Range("H137").Value = 1244.5714
$ws.Range("I137").Value = 1223.8125
$ws.Range("J137").Value = 1311
$ws.Range("K137").Value = 3671.4375
$ws.Range("L137").Value = 3933
$ws.Range("M137").Value = -1121.4375
$ws.Range("N137").Value = -9033

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21064.8
$ws.Range("I2").Value = 1103.6666
$ws.Range("J2").Value = 51006.5
$ws.Range("K2").Value = 1103.6666
$ws.Range("L2").Value = 51006.5
$ws.Range("M2").Value = -990.6666
$ws.Range("N2").Value = -51232.5

$ws.Range("H32").Value = 4519.8203
$ws.Range("I32").Value = 4236.543
$ws.Range("J32").Value = 6998.5
$ws.Range("K32").Value = 4236.543
$ws.Range("L32").Value = 6998.5
$ws.Range("M32").Value = -3949.543
$ws.Range("N32").Value = -7572.5

$ws.Range("H45").Value = 1965.1428
$ws.Range("I45").Value = 2208.4
$ws.Range("J45").Value = 1357
$ws.Range("K45").Value = 2208.4
$ws.Range("L45").Value = 1357
$ws.Range("M45").Value = -1831.4
$ws.Range("N45").Value = -2111

$ws.Range("H74").Value = 1017.13043
$ws.Range("I74").Value = 819
$ws.Range("K74").Value = 819
$ws.Range("M74").Value = 55

$ws.Range("H77").Value = 1017.13043
$ws.Range("I77").Value = 819
$ws.Range("K77").Value = 4095
$ws.Range("M77").Value = 273

$ws.Range("H116").Value = 21064.8
$ws.Range("I116").Value = 1103.6666
$ws.Range("J116").Value = 51006.5
$ws.Range("K116").Value = 1103.6666
$ws.Range("L116").Value = 51006.5
$ws.Range("M116").Value = 1190.3334
$ws.Range("N116").Value = -55594.5

$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -29820

$ws.Range("H132").Value = 2157.6365
$ws.Range("I132").Value = 1814.2106
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5442.6318
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -2912.6318
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21064.8
$ws.Range("I3").Value = 1103.6666
$ws.Range("J3").Value = 51006.5
$ws.Range("K3").Value = 1103.6666
$ws.Range("L3").Value = 51006.5
$ws.Range("M3").Value = -989.6666
$ws.Range("N3").Value = -51234.5

$ws.Range("H86").Value = 2582.818
$ws.Range("I86").Value = 2713.3635
$ws.Range("K86").Value = 2713.3635
$ws.Range("M86").Value = -1590.3635

$ws.Range("H89").Value = 2582.818
$ws.Range("I89").Value = 2713.3635
$ws.Range("K89").Value = 13566.8175
$ws.Range("M89").Value = -7950.817499999999

$ws.Range("H92").Value = 20399
$ws.Range("J92").Value = 20399
$ws.Range("L92").Value = 20399
$ws.Range("N92").Value = -25391

$ws.Range("H99").Value = 33334750
$ws.Range("I99").Value = 38462844
$ws.Range("K99").Value = 38462844
$ws.Range("M99").Value = -38461346

$ws.Range("H105").Value = 142859420
$ws.Range("I105").Value = 166669000
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 166669000
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = -166667253
$ws.Range("N105").Value = -5505

$ws.Range("H107").Value = 1438.85
$ws.Range("I107").Value = 1054
$ws.Range("J107").Value = 2978.25
$ws.Range("K107").Value = 1054
$ws.Range("L107").Value = 2978.25
$ws.Range("M107").Value = 866
$ws.Range("N107").Value = -6818.25

$ws.Range("H110").Value = 26975
$ws.Range("J110").Value = 26975
$ws.Range("L110").Value = 26975
$ws.Range("N110").Value = -35155

$ws.Range("H134").Value = 7088.2856
$ws.Range("I134").Value = 1520.8823
$ws.Range("J134").Value = 30749.75
$ws.Range("K134").Value = 4562.6469
$ws.Range("L134").Value = 92249.25
$ws.Range("M134").Value = -2027.6469
$ws.Range("N134").Value = -97319.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66667760
$ws.Range("I16").Value = 83334376
$ws.Range("J16").Value = 1304.3334
$ws.Range("K16").Value = 83334376
$ws.Range("L16").Value = 1304.3334
$ws.Range("M16").Value = -83334089
$ws.Range("N16").Value = -1878.3334

$ws.Range("H86").Value = 6082190.5
$ws.Range("J86").Value = 31889.143
$ws.Range("L86").Value = 31889.143
$ws.Range("N86").Value = -34135.143

$ws.Range("H89").Value = 6082190.5
$ws.Range("J89").Value = 31889.143
$ws.Range("L89").Value = 159445.715
$ws.Range("N89").Value = -170677.715

$ws.Range("H99").Value = 1810.7778
$ws.Range("I99").Value = 1779.1818
$ws.Range("K99").Value = 1779.1818
$ws.Range("M99").Value = -281.1818000000001

$ws.Range("H107").Value = 522.871
$ws.Range("I107").Value = 398.6
$ws.Range("J107").Value = 748.8182
$ws.Range("K107").Value = 398.6
$ws.Range("L107").Value = 748.8182
$ws.Range("M107").Value = 1521.4
$ws.Range("N107").Value = -4588.8182

$ws.Range("H109").Value = 24250.125
$ws.Range("J109").Value = 24250.125
$ws.Range("L109").Value = 24250.125
$ws.Range("N109").Value = -26330.125

$ws.Range("H113").Value = 66667760
$ws.Range("I113").Value = 83334376
$ws.Range("J113").Value = 1304.3334
$ws.Range("K113").Value = 83334376
$ws.Range("L113").Value = 1304.3334
$ws.Range("M113").Value = -83332206
$ws.Range("N113").Value = -5644.3334

$ws.Range("H126").Value = 1810.7778
$ws.Range("I126").Value = 1779.1818
$ws.Range("K126").Value = 5337.5454
$ws.Range("M126").Value = -2867.5454

$ws.Range("H132").Value = 4013.625
$ws.Range("I132").Value = 3542.2
$ws.Range("J132").Value = 4799.3335
$ws.Range("K132").Value = 10626.6
$ws.Range("L132").Value = 14398.0005
$ws.Range("M132").Value = -8096.599999999999
$ws.Range("N132").Value = -19458.0005

$ws.Range("H134").Value = 13889736
$ws.Range("I134").Value = 15873732
$ws.Range("J134").Value = 1766.6666
$ws.Range("K134").Value = 47621196
$ws.Range("L134").Value = 5299.9998
$ws.Range("M134").Value = -47618661
$ws.Range("N134").Value = -10369.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 395
$ws.Range("I23").Value = 480
$ws.Range("J23").Value = 310
$ws.Range("K23").Value = 1440
$ws.Range("L23").Value = 930
$ws.Range("M23").Value = -1205
$ws.Range("N23").Value = -1400

$ws.Range("H38").Value = 133.66667
$ws.Range("I38").Value = 100.5
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 301.5
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = 45.5
$ws.Range("N38").Value = -1294

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = $null

$ws.Range("H107").Value = 876.7143
$ws.Range("J107").Value = 876.7143
$ws.Range("L107").Value = 2630.1429
$ws.Range("N107").Value = -6470.1429

$ws.Range("H118").Value = 500
$ws.Range("I118").Value = 500
$ws.Range("K118").Value = 1500
$ws.Range("M118").Value = -257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2034.3077
$ws.Range("I102").Value = 2231
$ws.Range("K102").Value = 2231
$ws.Range("M102").Value = -609

$ws.Range("H132").Value = 2105.125
$ws.Range("I132").Value = 1691.7142
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5075.142599999999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2545.142599999999
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2279.5715
$ws.Range("J7").Value = 2242.3333
$ws.Range("L7").Value = 2242.3333
$ws.Range("N7").Value = -2466.3333

$ws.Range("H16").Value = 1910.3334
$ws.Range("I16").Value = 1932.4
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 1932.4
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -1762.4
$ws.Range("N16").Value = -2140

$ws.Range("H122").Value = 35717680
$ws.Range("I122").Value = 83337260
$ws.Range("J122").Value = 2991.75
$ws.Range("K122").Value = 250011780
$ws.Range("L122").Value = 8975.25
$ws.Range("M122").Value = -250009330
$ws.Range("N122").Value = -13875.25

$ws.Range("H126").Value = 2279.5715
$ws.Range("J126").Value = 2242.3333
$ws.Range("L126").Value = 6726.999899999999
$ws.Range("N126").Value = -11666.9999

$ws.Range("H132").Value = 41678.56
$ws.Range("I132").Value = 1097.4
$ws.Range("J132").Value = 102550.3
$ws.Range("K132").Value = 3292.2
$ws.Range("L132").Value = 307650.9
$ws.Range("M132").Value = -762.2000000000003
$ws.Range("N132").Value = -312710.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17336402
$ws.Range("I122").Value = 17336402
$ws.Range("K122").Value = 52009206
$ws.Range("M122").Value = -52006756

$ws.Range("H132").Value = 3559.5217
$ws.Range("I132").Value = 3243.389
$ws.Range("J132").Value = 4697.6
$ws.Range("K132").Value = 9730.167000000001
$ws.Range("L132").Value = 14092.8
$ws.Range("M132").Value = -7200.167000000001
$ws.Range("N132").Value = -19152.8

